$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "script shell -> Clement" -> "script shell Clement"  (drop the arrow)
# ---------------------------------------------------------------------------
$null = $d.Content.Find.Execute(
    "script shell -> Clément",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "script shell Clément",
    2)

# ---------------------------------------------------------------------------
# 2) Merge the split runs ("Groupe compose...") into a single run (no visible
#    text change, just collapses the spell-check run breaks).
# ---------------------------------------------------------------------------
$pGroupe = $d.Paragraphs.Item(3)
$null = $pGroupe.Range.Find.Execute(
    "Groupe composé de Clément Le Coadou et de Pauline Maceiras, groupe 1 MI Cy Tech",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Groupe composé de Clément Le Coadou et de Pauline Maceiras, groupe 1 MI Cy Tech",
    2)

# ---------------------------------------------------------------------------
# 3) Merge the split runs in "La correction des bugs..." into one run.
# ---------------------------------------------------------------------------
$pBugsOld = $d.Paragraphs.Item(7)
$null = $pBugsOld.Range.Find.Execute(
    "La correction des bugs liés à la mise en commun de nos programmes ainsi que les test se sont déroulés dans la dernière semaine avant la remise du projet.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "La correction des bugs liés à la mise en commun de nos programmes ainsi que les test se sont déroulés dans la dernière semaine avant la remise du projet.",
    2)

# ---------------------------------------------------------------------------
# 4) Replace the "…" placeholder right after "Partie langage c, fonctionnalités
#    manquantes : " with the actual missing-features text.
# ---------------------------------------------------------------------------
$pLangageC = $d.Paragraphs.Item(10)
$pLangageC.Range.Text = "L’option -r, -t3, -p3 et la gestion des erreurs n’est pas réalisé par le fichier c. Le graphique de sortie de -w ne semble pas conforme en vue de la taille des vecteurs et de leur orientation et -m ne fonctionne pas correctement."

# ---------------------------------------------------------------------------
# 5) Merge the runs in "Partie script shell, fonctionnalités manquantes : "
#    and drop the narrow no-break space before the colon.
# ---------------------------------------------------------------------------
$pScriptShellHeading = $d.Paragraphs.Item(11)
$null = $pScriptShellHeading.Range.Find.Execute(
    "Partie script shell, fonctionnalités manquantes : ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Partie script shell, fonctionnalités manquantes : ",
    2)

# ---------------------------------------------------------------------------
# 6) Merge the runs in "Il y a une fonctionnalité manquante..." and drop the
#    narrow no-break space before the colon.
# ---------------------------------------------------------------------------
$pMissingD = $d.Paragraphs.Item(12)
$null = $pMissingD.Range.Find.Execute(
    "Il y a une fonctionnalité manquante dans le script shell : le -d qui devait permettre à l’utilisateur de choisir une date min et max. ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Il y a une fonctionnalité manquante dans le script shell : le -d qui devait permettre à l’utilisateur de choisir une date min et max. ",
    2)

# ---------------------------------------------------------------------------
# 7) Merge the runs in "Le code pour réaliser..." and drop the narrow
#    no-break spaces around "????- ??- ??".
# ---------------------------------------------------------------------------
$pWhileShift = $d.Paragraphs.Item(13)
$null = $pWhileShift.Range.Find.Execute(
    "Le code pour réaliser cette fonction aurait dû utiliser un while/shift (comme pour le -f) pour récupérer la date min et max (avec un case par la suite de la forme ????- ??- ??) et la commande awk (similaire au tri de la localisation) pour trier les dates.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Le code pour réaliser cette fonction aurait dû utiliser un while/shift (comme pour le -f) pour récupérer la date min et max (avec un case par la suite de la forme ????- ??- ??) et la commande awk (similaire au tri de la localisation) pour trier les dates.",
    2)

# ---------------------------------------------------------------------------
# 8) Merge the runs in "Les options t3 et p3..." into one run.
# ---------------------------------------------------------------------------
$pT3P3 = $d.Paragraphs.Item(15)
$null = $pT3P3.Range.Find.Execute(
    "Les options t3 et p3 n’étant pas traitées par le langage c, il paraissait donc logique de ne pas les faire apparaître dans le script shell.",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Les options t3 et p3 n’étant pas traitées par le langage c, il paraissait donc logique de ne pas les faire apparaître dans le script shell.",
    2)

# ---------------------------------------------------------------------------
# 9) Merge the runs in "Partie Gnuplot : " and drop the narrow no-break
#    space before the colon.
# ---------------------------------------------------------------------------
$pGnuplotHeading = $d.Paragraphs.Item(16)
$null = $pGnuplotHeading.Range.Find.Execute(
    "Partie Gnuplot : ",
    $false, $false, $false, $false, $false,
    $true, 1, $false,
    "Partie Gnuplot : ",
    2)

# ---------------------------------------------------------------------------
# 10) Replace the "…" placeholder right after "Partie Gnuplot : " with the
#     actual missing-features text.
# ---------------------------------------------------------------------------
$pGnuplot = $d.Paragraphs.Item(17)
$pGnuplot.Range.Text = "Pour les options t2 et p2 le graphique de sortie n’est pas très lisible sûrement à cause des valeurs en x trop nombreuses. Une solution n’a pas était trouvée. "

# ---------------------------------------------------------------------------
# 11) Insert the new "L'écriture des scripts gnuplot..." paragraph right
#     after "La correction des bugs ... avant la remise du projet." and drop
#     the now-redundant blank paragraph that used to sit there.
# ---------------------------------------------------------------------------
$pBugs = $d.Paragraphs.Item(7)
$pBugs.Range.InsertParagraphAfter()
$pNew = $d.Paragraphs.Item(8)
$pNew.Range.Text = "L’écriture des scripts gnuplot a été réalisée en commun."

$pBlank = $d.Paragraphs.Item(9)
$pBlank.Range.Delete()
